# Append a new "05-Jan-2019" daily snapshot row to every metrics sheet in
# the workbook (8 sheets), mirroring the latest existing row's figures
# (the Coverity sheets carry a slightly updated "look-back" total).
#
# The date label must land as literal text (a shared string), exactly like
# the existing "DD-MMM-YYYY" rows already in the sheet, not as a real Excel
# date serial number. Entering it as a plain string gets auto-parsed into a
# date by the smart-typing layer, so instead we write it as a formula that
# evaluates to the text, then collapse that formula down to its plain
# cached value with Copy / PasteSpecial (values only) -- matching how the
# existing "03-Jan-2019" / "04-Jan-2019" cells are already stored (no
# number-format / quote-prefix side effects on the style table).

$wb = $excel.ActiveWorkbook

$xlPasteValues = -4163

function Add-DateRow {
    param($SheetIndex, $RowNumber, $DateText, $Values)

    $ws = $wb.Worksheets.Item($SheetIndex)

    $dateCell = $ws.Cells.Item($RowNumber, 1)
    $dateCell.Formula = '="' + $DateText + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial($xlPasteValues)

    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($RowNumber, 2 + $i).Value = $Values[$i]
    }
}

# 1: Wrong warning level -- row 5, all zero
Add-DateRow 1 5 "05-Jan-2019" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# 2: Treat warning not as error -- row 5, all zero
Add-DateRow 2 5 "05-Jan-2019" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# 3: Suppressed warnings -- row 5, all zero
Add-DateRow 3 5 "05-Jan-2019" @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

# 4: Actual warnings -- row 5, same totals as the prior three days
Add-DateRow 4 5 "05-Jan-2019" @(3336, 0, 420, 313, 378, 0, 133, 24, 1390, 24, 549, 33, 6600)

# 5: Coverity level 1 -- row 63, look-back window shrank (L: 15 -> 12)
Add-DateRow 5 63 "05-Jan-2019" @(9, 0, 8, 16, 0, 9, 17, 2, 51, 30, 12, 33, 187)

# 6: Coverity level 2 -- row 25, look-back window shrank (L: 539 -> 537)
Add-DateRow 6 25 "05-Jan-2019" @(461, 54, 127, 66, 62, 23, 290, 17, 282, 163, 537, 14, 2096)

# 7: Security level 1 -- row 5, same totals as 04-Jan-2019
Add-DateRow 7 5 "05-Jan-2019" @(6, 3, 0, 6, 23, 92, 277, 16, 83, 214, 297, 4, 1021)

# 8: Security level 2 -- row 5, same totals as 04-Jan-2019
Add-DateRow 8 5 "05-Jan-2019" @(38, 156, 2, 0, 0, 0, 319, 0, 1, 2, 0, 0, 518)
